$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 30
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 1
